# Add a new "TimeZone" column (column G) to the participant data sheet,
# filling every participant row with "Europe/Amsterdam".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("G1").Value = "TimeZone"

# Data rows (participants are in rows 2-13)
$lastRow = $ws.Cells.Item(1048576, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 13 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = "Europe/Amsterdam"
}

# Size the new column the way Excel does after typing values in and
# letting it best-fit to the widest entry ("Europe/Amsterdam").
$ws.Columns.Item(7).ColumnWidth = 9

# Move the active selection, mirroring the author's saved cursor position.
$ws.Range("F17").Select() | Out-Null
